# Apply the "Uploading newest EPS-US files" update:
#  - About sheet title becomes plural ("Exponents")
#  - "ETLE" sheet is renamed "ETLE-output" (exponent value -3 -> -4)
#  - A new "ETLE-capacity" sheet is added right after it (copy of ETLE-output
#    with exponent value -90), and becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- About sheet: pluralize the title -------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("A1").Value = "ETLE Electricity Technology Logit Exponents"

# --- Rename ETLE -> ETLE-output and update its exponent value -------------
$etleOutput = $wb.Worksheets.Item("ETLE")
$etleOutput.Name = "ETLE-output"
$etleOutput.Range("B2").Value = -4

# --- Duplicate it into a new "ETLE-capacity" sheet right after ------------
$etleOutput.Copy($null, $etleOutput) | Out-Null
$etleCapacity = $wb.Worksheets.Item($etleOutput.Index + 1)
$etleCapacity.Name = "ETLE-capacity"
$etleCapacity.Range("B2").Value = -90

# --- Make the new sheet the active / selected one --------------------------
$etleCapacity.Select() | Out-Null
$etleCapacity.Range("B2").Select() | Out-Null
